$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (single decimal point) need to be
# forced to Text format first, so Excel keeps them as text (matching the source data)
# instead of silently converting them to numeric cells.
$ws.Range('D2').Value = '60.960.59'
$ws.Range('E2').Value = '  -3.58%  '
$ws.Range('D3').Value = '2.453.43'
$ws.Range('E3').Value = '  -5.49%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '548.70'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.22%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.42'
$ws.Range('D6').NumberFormat = "General"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.58%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -2.20%  '
$ws.Range('D9').Value = '2.450.88'
$ws.Range('E9').Value = '  -5.51%  '
$ws.Range('E10').Value = '  -6.04%  '
$ws.Range('E11').Value = '  -1.37%  '
$ws.Range('E12').Value = '  -5.91%  '
$ws.Range('E13').Value = '  -5.33%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '25.92'
$ws.Range('D14').NumberFormat = "General"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -6.31%  '
$ws.Range('D15').Value = '2.892.73'
$ws.Range('E15').Value = '  -5.29%  '
$ws.Range('E16').Value = '  -6.35%  '
$ws.Range('D17').Value = '60.865.26'
$ws.Range('E17').Value = '  -3.52%  '
$ws.Range('D18').Value = '2.445.95'
$ws.Range('E18').Value = '  -6.32%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.04'
$ws.Range('D19').NumberFormat = "General"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -6.61%  '
$ws.Range('E20').Value = '  -6.55%  '
$ws.Range('E21').Value = '  -5.99%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '318.69'
$ws.Range('D22').NumberFormat = "General"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -5.50%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.43'
$ws.Range('D24').NumberFormat = "General"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.86%  '
$ws.Range('B25').Value = 'SuiNetwork'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.77'
$ws.Range('D25').NumberFormat = "General"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').Value = '0.0₃0974'
$ws.Range('E26').Value = '  -7.14%  '
$ws.Range('D27').Value = '2.570.60'
$ws.Range('E27').Value = '  -5.72%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '541.63'
$ws.Range('D28').NumberFormat = "General"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -5.74%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').NumberFormat = "General"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.49'
$ws.Range('D30').NumberFormat = "General"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.09%  '
$ws.Range('E31').Value = '  -7.35%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.57'
$ws.Range('D32').NumberFormat = "General"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.20%  '
$ws.Range('E33').Value = '  -5.93%  '
$ws.Range('E34').Value = '  -6.35%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.56'
$ws.Range('D35').NumberFormat = "General"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -7.66%  '
$ws.Range('E36').Value = '  -9.43%  '
$ws.Range('E37').Value = '  -8.20%  '
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.376'
$ws.Range('D39').NumberFormat = "General"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -5.06%  '
$ws.Range('E40').Value = '  -5.12%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '145.16'
$ws.Range('D41').NumberFormat = "General"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -5.74%  '
$ws.Range('E42').Value = '  -6.48%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '39.88'
$ws.Range('D44').NumberFormat = "General"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.50%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.30'
$ws.Range('D45').NumberFormat = "General"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -6.37%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '145.58'
$ws.Range('D46').NumberFormat = "General"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -7.03%  '
$ws.Range('E47').Value = '  -6.59%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '20.78'
$ws.Range('D48').NumberFormat = "General"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -9.08%  '
$ws.Range('E49').Value = '  -7.88%  '
$ws.Range('E50').Value = '  -6.42%  '
$ws.Range('E51').Value = '  -4.57%  '

Write-Output "Applied cryptos update"
